$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 270, shifting all existing rows
# from 270 downward to 271 downward (new dimension becomes A1:R336).
$ws.Rows.Item(270).Insert()

# Populate the newly inserted row 270 with the new record.
$ws.Range("A270").Value = 4
$ws.Range("B270").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value = "Los Lagos"
$ws.Range("D270").Value = 44855
$ws.Range("E270").Value = 10
$ws.Range("F270").Value = 100112043
$ws.Range("G270").Value = "Pepino ensalada"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 400
$ws.Range("K270").Value = 25000
$ws.Range("L270").Value = 25000
$ws.Range("M270").Value = 25000
$ws.Range("N270").Value = "$/caja 60 unidades"
$ws.Range("O270").Value = "Región de Arica y Parinacota"
$ws.Range("P270").Value = 417
$ws.Range("Q270").Value = 60
$ws.Range("R270").Value = "Hortaliza"

# Preserve the date-number formatting of column D for the new row,
# matching the style used by the rest of the date column.
$ws.Range("D270").NumberFormat = $ws.Range("D271").NumberFormat
